# Weekly update: a new "Fruta" (Piña) price record was reported for
# Terminal Hortofrutícola Agro Chillán, dated 2023-03-28 (serial 45013).
# It belongs chronologically right before the existing row 275 record
# (2022-07-07 / serial 44568), so insert a new row at 275, push the rest
# of the table down by one, and populate the new row with the reported
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 275; everything from the old row
# 275 downward (through the old last row, 296) shifts down to 276-297,
# and the sheet dimension grows from A1:T296 to A1:T297 automatically.
$ws.Rows(275).Insert()

# Fill in the newly inserted row 275 with the new observation.
$ws.Range("A275").Value = 7
$ws.Range("B275").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C275").Value = "Ñuble"
$ws.Range("D275").Value = 45013
$ws.Range("E275").Value = 16
$ws.Range("F275").Value = "Fruta"
$ws.Range("G275").Value = 100108
$ws.Range("H275").Value = "Tropicales y subtropicales"
$ws.Range("I275").Value = 100108005
$ws.Range("J275").Value = "Piña"
$ws.Range("K275").Value = "Caramelo"
$ws.Range("L275").Value = "Segunda"
$ws.Range("M275").Value = 30
$ws.Range("N275").Value = 25000
$ws.Range("O275").Value = 25000
$ws.Range("P275").Value = 25000
$ws.Range("Q275").Value = "$/caja 14 unidades"
$ws.Range("R275").Value = "Ecuador"
$ws.Range("S275").Value = 1786
$ws.Range("T275").Value = 14
